# Table1.xlsx re-upload: the top of the sheet had a merged "Heading"
# banner row (A1:G1, shared string "Heading") sitting above the real
# header row. That banner row is removed entirely, so the column
# headers ("Sr. No.", "Energy applied (J) to UUC", ...) become row 1
# and the six data rows shift up to rows 2-7.
#
# Deleting the whole row (rather than just clearing it) is what makes
# Excel drop the now-unused "Heading" shared string, the A1:G1 merged
# cell, and the style/border combo that only that banner row used.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Rows.Item(1).Delete()

# Matches the saved selection in the edited file.
$ws.Range("F12").Select()
